$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("G2").Value = 10.08846466666667
$ws.Range("H2").Value = 30.265394
$ws.Range("I2").Value = 0.2597171077778241
$ws.Range("J2").Value = 0.2597171077778241
$ws.Range("K2").Value = 2
$ws.Range("L2").Value = 0.6666666666666666
$ws.Range("M2").Value = 0.2773296666666667
$ws.Range("N2").Value = 0.8319890000000001
$ws.Range("Q2").Value = 2.797830543185111
$ws.Range("R2").Value = 25.180474888666
$ws.Range("S2").Value = 0.2597171077778241
$ws.Range("T2").Value = 0.2597171077778241

# Row 3
$ws.Range("I3").Value = 0.05876531725312483
$ws.Range("J3").Value = 0.05876531725312482
$ws.Range("K3").Value = 2
$ws.Range("L3").Value = 0.6666666666666666
$ws.Range("M3").Value = 0.2773296666666667
$ws.Range("N3").Value = 0.8319890000000001
$ws.Range("Q3").Value = 0.6330557154956667
$ws.Range("R3").Value = 5.697501439461
$ws.Range("S3").Value = 0.05876531725312483
$ws.Range("T3").Value = 0.05876531725312482

# Row 4
$ws.Range("G4").Value = 14.90894133333333
$ws.Range("H4").Value = 44.726824
$ws.Range("I4").Value = 0.3838153030278664
$ws.Range("J4").Value = 0.3838153030278664
$ws.Range("K4").Value = 2
$ws.Range("L4").Value = 0.6666666666666666
$ws.Range("M4").Value = 0.2773296666666667
$ws.Range("N4").Value = 0.8319890000000001
$ws.Range("Q4").Value = 4.134691730326223
$ws.Range("R4").Value = 37.21222557293601
$ws.Range("S4").Value = 0.3838153030278664
$ws.Range("T4").Value = 0.3838153030278664

# Row 5
$ws.Range("E5").Value = 3
$ws.Range("F5").Value = 1
$ws.Range("G5").Value = 0.2836386666666667
$ws.Range("H5").Value = 0.850916
$ws.Range("I5").Value = 0.007301984652235982
$ws.Range("J5").Value = 0.007301984652235982
$ws.Range("K5").Value = 2
$ws.Range("L5").Value = 0.6666666666666666
$ws.Range("M5").Value = 0.2773296666666667
$ws.Range("N5").Value = 0.8319890000000001
$ws.Range("Q5").Value = 0.07866141688044445
$ws.Range("R5").Value = 0.7079527519240001
$ws.Range("S5").Value = 0.007301984652235982
$ws.Range("T5").Value = 0.007301984652235982

# Row 6
$ws.Range("G6").Value = 11.28032366666667
$ws.Range("H6").Value = 33.840971
$ws.Range("I6").Value = 0.2904002872889486
$ws.Range("J6").Value = 0.2904002872889486
$ws.Range("K6").Value = 2
$ws.Range("L6").Value = 0.6666666666666666
$ws.Range("M6").Value = 0.2773296666666667
$ws.Range("N6").Value = 0.8319890000000001
$ws.Range("Q6").Value = 3.128368402368778
$ws.Range("R6").Value = 28.155315621319
$ws.Range("S6").Value = 0.2904002872889486
$ws.Range("T6").Value = 0.2904002872889486

